# Applies the "Complete all figure captions" commit to the lab report.
#
# Strategy: each "Fig. N.  Caption" paragraph has its trailing "Caption"
# word deleted (via Find, scoped to that single paragraph so the other
# identical placeholders are untouched) and then the real caption text is
# typed back in as a sequence of separate Insert-After calls - this mirrors
# how Word COM naturally produces one new <w:r> per discrete insertion,
# matching the run layout seen in the target document.

$d = $word.ActiveDocument

function Remove-CaptionPlaceholder($paraIndex) {
    $p = $d.Paragraphs($paraIndex)
    $r = $p.Range
    $r.Find.ClearFormatting()
    [void]$r.Find.Execute("Caption", $false, $false, $false, $false, $false, $true, 1, $false, "", 2)
}

function Add-CaptionRun($paraIndex, [string]$text) {
    $p = $d.Paragraphs($paraIndex)
    $r = $p.Range
    $r.MoveEnd(1, -1) | Out-Null
    $r.Collapse(0)
    $r.InsertAfter($text)
}

# ---------------------------------------------------------------------
# Fig. 3 -- "Plot of channel data (scrolling) ... while moving my left hand."
# ---------------------------------------------------------------------
Remove-CaptionPlaceholder 41
Add-CaptionRun 41 "Plot of channel data (scrolling), represented in the time domain, while "
Add-CaptionRun 41 "moving "
Add-CaptionRun 41 "my"
Add-CaptionRun 41 " left hand"
Add-CaptionRun 41 "."

# ---------------------------------------------------------------------
# Fig. 4 -- frequency-domain counterpart, plus extra discussion sentences.
# ---------------------------------------------------------------------
Remove-CaptionPlaceholder 43
Add-CaptionRun 43 "Plot of channel spectra and maps, represented in the frequency domain, while "
Add-CaptionRun 43 "moving "
Add-CaptionRun 43 "my"
Add-CaptionRun 43 " left hand"
Add-CaptionRun 43 "."
Add-CaptionRun 43 " Note the prevalence of neural activity in the right cortical hemisphere. This instance of contralateral activation supports the recent studies mentioned in the abstract."

# ---------------------------------------------------------------------
# Fig. 5 -- imagining moving left hand, time domain.
# ---------------------------------------------------------------------
Remove-CaptionPlaceholder 45
Add-CaptionRun 45 "Plot of channel data (scrolling), represented in the time domain, while "
Add-CaptionRun 45 "imagining "
Add-CaptionRun 45 "moving "
Add-CaptionRun 45 "my"
Add-CaptionRun 45 " left hand."

# ---------------------------------------------------------------------
# Fig. 6 -- imagining moving left hand, frequency domain, plus discussion.
# ---------------------------------------------------------------------
Remove-CaptionPlaceholder 47
Add-CaptionRun 47 "Plot of channel spectra and maps, represented in the frequency domain, while "
Add-CaptionRun 47 "imagining "
Add-CaptionRun 47 "moving "
Add-CaptionRun 47 "my"
Add-CaptionRun 47 " left hand."
Add-CaptionRun 47 " Note that, although there is a presence of contralateral neural activity, it is hardly comparable to magnitudes observed while actually moving the hand (as seen in Fig. 4)."

# ---------------------------------------------------------------------
# Fig. 7 -- moving right hand, time domain, plus two discussion sentences.
# ---------------------------------------------------------------------
Remove-CaptionPlaceholder 49
Add-CaptionRun 49 "Plot of channel data (scrolling), represented in the time domain, while moving "
Add-CaptionRun 49 "my"
Add-CaptionRun 49 " "
Add-CaptionRun 49 "right"
Add-CaptionRun 49 " hand."
Add-CaptionRun 49 " Although visibly more erratic, note that the scaling factor here is drastically lower than those in other trials."
Add-CaptionRun 49 " It is hard to say why this is. Perhaps moving my right hand required less neural activity since I am right-handed, or I had not performed these motions with as much energy as I had while moving my left hand (Fig. 3-4)."

# ---------------------------------------------------------------------
# Fig. 8 -- moving right hand, frequency domain.
# ---------------------------------------------------------------------
Remove-CaptionPlaceholder 51
Add-CaptionRun 51 "Plot of channel spectra and maps, represented in the frequency domain, while moving "
Add-CaptionRun 51 "my"
Add-CaptionRun 51 " "
Add-CaptionRun 51 "right"
Add-CaptionRun 51 " hand."

# ---------------------------------------------------------------------
# Fig. 9 -- imagining moving right hand, time domain.
# ---------------------------------------------------------------------
Remove-CaptionPlaceholder 53
Add-CaptionRun 53 "P"
Add-CaptionRun 53 "lot of channel data (scrolling) represented in the time domain,"
Add-CaptionRun 53 " while "
Add-CaptionRun 53 "imagining "
Add-CaptionRun 53 "m"
Add-CaptionRun 53 "oving "
Add-CaptionRun 53 "my "
Add-CaptionRun 53 "right"
Add-CaptionRun 53 " hand"
Add-CaptionRun 53 "."

# ---------------------------------------------------------------------
# Fig. 10 -- imagining moving right hand, frequency domain, plus discussion,
# and the stray empty ListParagraph right after it is removed.
# ---------------------------------------------------------------------
Remove-CaptionPlaceholder 55
Add-CaptionRun 55 "Plot of channel spectra and maps, represented in the frequency domain, while "
Add-CaptionRun 55 "imagining "
Add-CaptionRun 55 "moving the "
Add-CaptionRun 55 "right"
Add-CaptionRun 55 " hand."
Add-CaptionRun 55 " Note that there is clearly an increase in neural activity in the left cortical hemisphere, which is again in accordance with the studies mentioned in the abstract."
Add-CaptionRun 55 " However, there is still some peculiar activity in the right caudal region of my brain, noticeable in 9-11 Hz "
Add-CaptionRun 55 "band of "
Add-CaptionRun 55 "the"
Add-CaptionRun 55 " μ rhythm."

# Remove the now-empty ListParagraph that used to sit between Fig. 10's
# caption and the "Discussion" heading.
$empty = $d.Paragraphs(56)
$empty.Range.Delete()

# ---------------------------------------------------------------------
# Move the _GoBack bookmark from the end of Fig. 2's caption to the start
# of the "Discussion" heading, and drop the stray leading space run there.
# ---------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

$discussionPara = $d.Paragraphs(56)
$discussionRange = $discussionPara.Range
[void]$discussionRange.Find.Execute(" Discussion", $false, $false, $false, $false, $false, $true, 1, $false, "Discussion", 2)

$discussionPara2 = $d.Paragraphs(56)
$startRange = $discussionPara2.Range
$startRange.Collapse(1)
$d.Bookmarks.Add("_GoBack", $startRange)
